# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    shows it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) Shrink the (now narrower) status columns:
#      Overview columns E & F, zh-cn column C, de-de column C
#    from ~17.22 chars down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Column width achievable through Excel's pixel-quantized ColumnWidth
# property that lands closest to the target stored width (13.4101845877511).
$newColWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update status values -------------------------------------------------
if ($overview.Range("E2").Value() -eq $oldStatus) {
    $overview.Range("E2").Value = $newStatus
}
if ($overview.Range("F2").Value() -eq $oldStatus) {
    $overview.Range("F2").Value = $newStatus
}
if ($zhcn.Range("C2").Value() -eq $oldStatus) {
    $zhcn.Range("C2").Value = $newStatus
}
if ($dede.Range("C2").Value() -eq $oldStatus) {
    $dede.Range("C2").Value = $newStatus
}

# --- Narrow the status columns ---------------------------------------------
$overview.Range("E1").ColumnWidth = $newColWidth
$overview.Range("F1").ColumnWidth = $newColWidth
$zhcn.Range("C1").ColumnWidth = $newColWidth
$dede.Range("C1").ColumnWidth = $newColWidth
